# Apply the "赣州·十万伏特-次元音乐only" insertion (2024-04-20) and the associated
# "想去人数" (F column) refresh to the 江西-漫展信息 workbook.
#
# Two worksheets are affected:
#   展览   (sheet index 1) - insert new row at row 8
#   全部类型 (sheet index 4) - insert new row at row 9 (it has one extra lead-in
#                              row, "南昌·Kpop New Life", that 展览 does not have)
#
# In both sheets, inserting the row pushes every following row down by one, and a
# handful of the "想去人数" (F column) figures were also refreshed upward as part of
# the same re-scrape that added the new row.

$wb = $excel.ActiveWorkbook

function Update-Sheet($ws, [int]$insertRow, [int]$totalRowsAfterInsert, [int]$newWorldRow) {
    # 1) Refresh the "New World国潮动漫博览会" attendance figure in place (no shift
    #    involved for this row, it sits above the insertion point).
    $ws.Cells.Item($newWorldRow, 6).Value = 4842

    # 2) Insert a new blank row, shifting row $insertRow and everything below it down
    #    by one.
    $ws.Rows.Item($insertRow).Insert()

    # 3) Restore the column-A number formatting (bold + border) on the new row by
    #    copying the formatting only from the row right below it (which now holds
    #    what used to be row $insertRow, and already carries the correct style).
    $ws.Range("A" + ($insertRow + 1)).Copy() | Out-Null
    $ws.Range("A" + $insertRow).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0

    # 4) Populate the newly inserted row with the "赣州·十万伏特-次元音乐only" event.
    #    Column B holds a plain text date string ("2024-04-20"); Excel's COM layer
    #    would otherwise auto-convert that into a real date value, so force a text
    #    number format while writing it, then restore "General" to match the rest
    #    of the column.
    $bCell = $ws.Cells.Item($insertRow, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = "2024-04-20"
    $bCell.NumberFormat = "General"
    $ws.Cells.Item($insertRow, 3).Value = "赣州·十万伏特-次元音乐only"
    $ws.Cells.Item($insertRow, 4).Value = "平安大道 麋鹿LiveHouse"
    $ws.Cells.Item($insertRow, 5).Value = "2024.04.20 14:30-04.21 21:00"
    $ws.Cells.Item($insertRow, 6).Value = 0
    $ws.Cells.Item($insertRow, 7).Value = 55
    $ws.Cells.Item($insertRow, 8).Value = "https://show.bilibili.com/platform/detail.html?id=83871"
    $ws.Cells.Item($insertRow, 9).Value = "//i2.hdslb.com/bfs/openplatform/202404/XlyMkr9u1711959548249.jpeg"

    # 5) A handful of "想去人数" values for the rows that got shifted down were
    #    refreshed too (same re-scrape). Offsets are expressed relative to
    #    $insertRow (i.e. the row number *after* the insert/shift).
    $offsets = 1, 2, 6, 8, 11, 12, 14, 16, 17
    $newValues = 534, 486, 3234, 118, 2490, 117, 34, 30, 119
    for ($i = 0; $i -lt $offsets.Length; $i++) {
        $row = $insertRow + $offsets[$i]
        $ws.Cells.Item($row, 6).Value = $newValues[$i]
    }

    # 6) The leading "序号" column (A) just stores row-number-minus-one for every
    #    data row; recompute it for all rows so it stays consistent after the
    #    insert.
    for ($r = 2; $r -le $totalRowsAfterInsert; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

$wsExhibit = $wb.Worksheets.Item("展览")
Update-Sheet $wsExhibit 8 28 5

$wsAll = $wb.Worksheets.Item("全部类型")
Update-Sheet $wsAll 9 29 6
